$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'paper',
    '101 hanasoge aa14345-10.pdf',
    '107 heister 2010-eccomas.pdf',
    '108 heister chp%3A10.1007%2F978-3-642-15646-5_13.pdf',
    '13 Geology-2013-Ballmer-G33804.1.pdf',
    '130 van summeren 2041-8205_736_1_L15.pdf',
    '131 kaneko 1-s2.0-S0040195110002684-main.pdf',
    '132 Geophys. J. Int.-2014-Kaneko-gji-ggu030.pdf',
    '134 key nature11932.pdf',
    '14 Ballmer ggge20662.pdf',
    '143 kreiss 110832847.pdf',
    '147 Geophys. J. Int.-2015-Lahivaara-1662-75.pdf',
    '148 lahivaara 0266-5611_30_1_015003.pdf',
    '150 lepourhiet ggge2265.pdf',
    '16 Barbot_et_al-2010-Geophysical_Journal_International.pdf',
    '164 ma Xiao_Ma.pdf',
    '168 manea Geology-2012-Manea-35-8.pdf',
    '17 Barbot_et_al-2010-Geophysical_Journal_International (1).pdf',
    '18 Barriere dbfe70dce0b8f0586cd083fe58cad33b5887f39483ffaa7c18a655e0f8b250f3.pdf',
    '183 michea 2c274b5da1f640f282336a91b4d159a812d792d17d8ffb62b1ebb9fc2bc78526.pdf',
    '184 michel 06236135.pdf',
    '189 motoki ggge20669.pdf',
    '191 Natarov_et_al-2012-Geophysical_Journal_International.pdf',
    '197 Obermaier_et_al-2011-Computer_Graphics_Forum.pdf',
    '203 ozgun 2621.full.pdf',
    '207 quere jgrb50102.pdf',
    '210 rappapport 0004-637X_752_1_1.pdf',
    '217 Geology-2012-Rey-e280.pdf',
    '219 Lithosphere-2010-Rey-328-32.pdf',
    '223 Geophys. J. Int.-2013-Ristau-gji-ggt222.pdf',
    '226 rollins art%3A10.1007%2Fs00024-014-1005-6.pdf',
    '227 rousset jgrb17375.pdf',
    '230 savage 1241.full.pdf',
    '25 bianco jgrb16872.pdf',
    '259 thebaud 1-s2.0-S0301926812000769-main.pdf',
    '26 billen1-s2.0-S0031920110001160-main.pdf',
    '268 trubitsyn art%3A10.1134%2FS1069351312060079.pdf',
    '269 trubitsyn art%3A10.1134%2FS106935131112010X.pdf',
    '276 vansumm 2041-8205_736_1_L15.pdf',
    '277 vansumm ggge2119.pdf',
    '278 vansumm jgre20077.pdf',
    '30 Boue Geophys. J. Int.-2013-Bouee-844-8.pdf',
    '307 Zhou_et_al-2011-Geophysical_Journal_International.pdf',
    '31 boue 1-s2.0-S0012821X14000636-main.pdf',
    '313 zhu grl29542.pdf',
    '35 ebru j.1365-246X.2009.04468.x.pdf',
    '36 ebru j.1365-246X.2011.04970.x.pdf',
    '37 Bruhat.pdf',
    '51 conrad Geological Society of America Bulletin-2013-Conrad-1027-52.pdf',
    '53 conrad nature12203.pdf',
    '7 arredondo 1-s2.0-S0031920112000349-main.pdf',
    '73 drilleau Dupros_parco_2010_ref6125.pdf',
    '74 Elkinsannurev-earth-042711-105503 (1).pdf',
    '85 Geology-2011-Flament-1159-62.pdf'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

$wb.Save()
